$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7282781
$ws.Range("C4").Value = 38597
$ws.Range("D4").Value = 4512774
$ws.Range("E4").Value = 2560897
$ws.Range("G4").Value = 670
$ws.Range("H4").Value = 209110

$ws.Range("D6").Value = 4050837
$ws.Range("E6").Value = 525748

$ws.Range("B8").Value = 806038
$ws.Range("C8").Value = 7721
$ws.Range("D8").Value = 700112
$ws.Range("E8").Value = 80630
$ws.Range("G8").Value = 193
$ws.Range("H8").Value = 25296

$ws.Range("B12").Value = 702484
$ws.Range("C12").Value = 11249
$ws.Range("D12").Value = 556489
$ws.Range("E12").Value = 130452
$ws.Range("G12").Value = 335
$ws.Range("H12").Value = 15543

$ws.Range("B13").Value = 669498
$ws.Range("C13").Value = 969
$ws.Range("D13").Value = 601818
$ws.Range("E13").Value = 51304
$ws.Range("G13").Value = 64
$ws.Range("H13").Value = 16376

$ws.Range("B37").Value = 110108
$ws.Range("C37").Value = 677
$ws.Range("D37").Value = 86796
$ws.Range("E37").Value = 20989
$ws.Range("G37").Value = 12
$ws.Range("H37").Value = 2323

$ws.Range("B41").Value = 102736
$ws.Range("C41").Value = 111
$ws.Range("D41").Value = 94374
$ws.Range("E41").Value = 2493
$ws.Range("G41").Value = 16
$ws.Range("H41").Value = 5869

$ws.Range("B48").Value = 81055
$ws.Range("C48").Value = 558
$ws.Range("D48").Value = 74151
$ws.Range("E48").Value = 5364
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 1540

$ws.Range("B58").Value = 58198
$ws.Range("C58").Value = 136
$ws.Range("D58").Value = 49722
$ws.Range("E58").Value = 7370
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 1106

$ws.Range("E71").Value = 12514
$ws.Range("G71").Value = 7
$ws.Range("H71").Value = 689

$ws.Range("D76").Value = 17832
$ws.Range("E76").Value = 13497

$ws.Range("B91").Value = 14612
$ws.Range("C91").Value = 97
$ws.Range("D91").Value = 13727
$ws.Range("E91").Value = 553

$ws.Range("A93").Value = "Noruega"
$ws.Range("B93").Value = 13627
$ws.Range("C93").Value = 82
$ws.Range("D93").Value = 11190
$ws.Range("E93").Value = 2167
$ws.Range("H93").Value = 270

$ws.Range("A94").Value = "Sudan"
$ws.Range("B94").Value = 13606
$ws.Range("C94").Value = 14
$ws.Range("D94").Value = 6764
$ws.Range("E94").Value = 6006
$ws.Range("H94").Value = 836

$ws.Range("B99").Value = 10512
$ws.Range("C99").Value = 34
$ws.Range("D99").Value = 9836
$ws.Range("E99").Value = 611

$ws.Range("A118").Value = "Suazilandia"
$ws.Range("B118").Value = 5419
$ws.Range("C118").Value = 20
$ws.Range("D118").Value = 4802
$ws.Range("E118").Value = 509
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 108

$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 5412
$ws.Range("C119").Value = 62
$ws.Range("D119").Value = 4732
$ws.Range("E119").Value = 560
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 120

$ws.Range("A120").Value = "Republica de Yibuti"
$ws.Range("B120").Value = 5409
$ws.Range("D120").Value = 5340
$ws.Range("E120").Value = 8
$ws.Range("H120").Value = 61

$ws.Range("B124").Value = 5028
$ws.Range("C124").Value = 10
$ws.Range("D124").Value = 4740
$ws.Range("E124").Value = 205

$ws.Range("A128").Value = "Ruanda"
$ws.Range("B128").Value = 4811
$ws.Range("C128").Value = 13
$ws.Range("D128").Value = 3091
$ws.Range("E128").Value = 1691
$ws.Range("H128").Value = 29

$ws.Range("A129").Value = "Republica de Africa Central"
$ws.Range("B129").Value = 4806
$ws.Range("D129").Value = 1840
$ws.Range("E129").Value = 2904
$ws.Range("H129").Value = 62

$ws.Range("A132").Value = "Trinidad yTobago"
$ws.Range("B132").Value = 4312
$ws.Range("C132").Value = 35
$ws.Range("D132").Value = 2185
$ws.Range("E132").Value = 2057
$ws.Range("H132").Value = 70

$ws.Range("A133").Value = "Lituania"
$ws.Range("B133").Value = 4295
$ws.Range("C133").Value = 111
$ws.Range("D133").Value = 2319
$ws.Range("E133").Value = 1887
$ws.Range("H133").Value = 89

$ws.Range("A138").Value = "Somalia"
$ws.Range("B138").Value = 3588
$ws.Range("C138").Value = 123
$ws.Range("D138").Value = 2943
$ws.Range("E138").Value = 546
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 99

$ws.Range("A139").Value = "Gambia"
$ws.Range("B139").Value = 3555
$ws.Range("D139").Value = 2034
$ws.Range("E139").Value = 1411
$ws.Range("H139").Value = 110

$ws.Range("A140").Value = "Mayotte"
$ws.Range("B140").Value = 3541
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 2964
$ws.Range("E140").Value = 537
$ws.Range("H140").Value = 40

$ws.Range("A141").Value = "Tailandia"
$ws.Range("B141").Value = 3522
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 3362
$ws.Range("E141").Value = 101
$ws.Range("H141").Value = 59

$ws.Range("B147").Value = 2725
$ws.Range("C147").Value = 16
$ws.Range("D147").Value = 1535
$ws.Range("E147").Value = 1116
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 74

$ws.Range("B177").Value = 485
$ws.Range("C177").Value = 2
$ws.Range("D177").Value = 472
$ws.Range("E177").Value = 12

$ws.Range("B178").Value = 478
$ws.Range("C178").Value = 4
$ws.Range("D178").Value = 458
$ws.Range("E178").Value = 13

$ws.Range("B179").Value = 460
$ws.Range("C179").Value = 2
$ws.Range("E179").Value = 43

$ws.Range("A206").Value = "Timor Oriental"

$ws.Range("A207").Value = "Santa Lucia"
